$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T_BLGD")

# Activate T_BLGD (this sheet becomes the active tab; DEFECTOS loses tabSelected)
$ws.Activate()

# Populate the new part-number rows. Values are written in the same
# left-to-right / top-to-bottom first-use order as the source data so new
# shared-string entries land at the expected indices.
$ws.Range("A14").Value = "RIGID.1-OMEGA"
$ws.Range("B14").Value = "L528-82711-00301A     "
$ws.Range("A15").Value = "RIGID. 2-OMEGA"
$ws.Range("B15").Value = "L528-82712-00301A     "
$ws.Range("A16").Value = "TAPA"
$ws.Range("B16").Value = "L528-82923-001     "
$ws.Range("A17").Value = "OMEGA 1"
$ws.Range("A18").Value = "OMEGA 2"
$ws.Range("B17").Value = "L528-82211-005   "
$ws.Range("B18").Value = "L528-82212-007    "

# Column C mirrors column A on each row. A14/A17 carry the bold "s=11"
# style, so copy their formatting across; the rest stay default-styled.
$ws.Range("A14").Copy($ws.Range("C14"))
$ws.Range("C15").Value = "RIGID. 2-OMEGA"
$ws.Range("C16").Value = "TAPA"
$ws.Range("A17").Copy($ws.Range("C17"))
$ws.Range("C18").Value = "OMEGA 2"

# Column B widened slightly to fit the new part numbers.
$ws.Columns.Item(2).ColumnWidth = 14.14

# Leave the selection on the last new row, matching the saved view state.
$ws.Range("A18").Select()
